# Week_7_Logs.xlsx -- "copied over local copies of logs"
#
# Jesse fills in his personal copies of the TASK SUMMARY SHEET and the
# ACTIVITY LOG SUMMARY SHEET for week 7 (name, week number, and the task
# rows / hour totals), and leaves the TASK SUMMARY SHEET as the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# TASK SUMMARY SHEET
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Header: who + which week this is for.
# ("Project Build" must be written before "Jesse Hare" so the shared
# string table indices line up with the authored workbook.)
$ws3.Range("A3").Value = "Project Build"
$ws3.Range("C1").Value = "Jesse Hare"
$ws3.Range("E1").Value = 7

# Task rows for the week.
$ws3.Range("B3").Value = "Continue work on dynamic search/filter"
$ws3.Range("C3").Value = 6
$ws3.Range("D3").Value = 6
$ws3.Range("E3").Value = 6

$ws3.Range("A4").Value = "Project Build"
$ws3.Range("B4").Value = "Fix issues with the Search function"
$ws3.Range("C4").Value = 2
$ws3.Range("D4").Value = 2
$ws3.Range("E4").Value = 2

$ws3.Range("A5").Value = "Project Build"
$ws3.Range("B5").Value = "Add tooltips to UI elements"
$ws3.Range("C5").Value = 1
$ws3.Range("D5").Value = 1
$ws3.Range("E5").Value = 0

$ws3.Range("A6").Value = "Project Build"
$ws3.Range("B6").Value = "Error handling and input validation"
$ws3.Range("C6").Value = 5
$ws3.Range("D6").Value = 7
$ws3.Range("E6").Value = 0

$ws3.Range("A7").Value = "Project Build"
$ws3.Range("B7").Value = "Query optimisation"
$ws3.Range("C7").Value = 1
$ws3.Range("D7").Value = 1
$ws3.Range("E7").Value = 0

$ws3.Range("A8").Value = "Project Build"
$ws3.Range("B8").Value = "Add sort by header fnctionality (click on header to sort)"
$ws3.Range("C8").Value = 3
$ws3.Range("D8").Value = 3
$ws3.Range("E8").Value = 0

# ---------------------------------------------------------------------
# ACTIVITY LOG SUMMARY SHEET
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("ACTIVITY LOG SUMMARY SHEET")

$ws4.Range("D1").Value = "Jesse Hare"
$ws4.Range("A4").Value = "Project Build"
$ws4.Range("B4").Value = 15
$ws4.Range("C4").Value = 5

# ---------------------------------------------------------------------
# Leave the TASK SUMMARY SHEET selected/active, matching the saved
# selection left behind in the workbook.
# ---------------------------------------------------------------------
[void]$ws3.Activate()
[void]$ws3.Range("B3:C4").Select()
